# Commit: "added harvard case classification"
# The underlying per-app evaluation stats (for the "_old" app variants) were
# recomputed after a new Harvard case was added to the dataset, which shifts the
# average/variance/std-dev figures for Ada_old, Avey_old, Babylon_old, Buoy_old,
# K health_old, WebMD_old, doctor_MA_old, doctor_NJ_old and doctor_TH_old, and the
# two summary columns (average_doctor / average_doctor_old) swap position with the
# now-recalculated values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: header swap for average_doctor / average_doctor_old
$ws.Range("BP1").Value = "average_doctor_old"
$ws.Range("BQ1").Value = "average_doctor"

# Row 4: recalculated stats
$ws.Range("AI4").Value = 0.26
$ws.Range("AJ4").Value = 0.11
$ws.Range("AK4").Value = 0.332
$ws.Range("AU4").Value = 0.279
$ws.Range("AV4").Value = 0.02
$ws.Range("AW4").Value = 0.143
$ws.Range("BA4").Value = 2.051
$ws.Range("BB4").Value = 0.08699999999999999
$ws.Range("BC4").Value = 0.294
$ws.Range("BG4").Value = 0.708
$ws.Range("BH4").Value = 0.179
$ws.Range("BI4").Value = 0.423
$ws.Range("BM4").Value = 0.708
$ws.Range("BN4").Value = 0.089
$ws.Range("BO4").Value = 0.298
$ws.Range("BP4").Value = 0.6840000000000001
$ws.Range("BQ4").Value = 0.8110000000000001
$ws.Range("E4").Value = 0.508
$ws.Range("F4").Value = 0.049
$ws.Range("G4").Value = 0.221
$ws.Range("N4").Value = 0.594
$ws.Range("O4").Value = 0.103
$ws.Range("P4").Value = 0.321
$ws.Range("Q4").Value = 0.275
$ws.Range("R4").Value = 0.124
$ws.Range("S4").Value = 0.352
$ws.Range("W4").Value = 0.396
$ws.Range("X4").Value = 0.118
$ws.Range("Y4").Value = 0.343

# Row 5: recalculated stats
$ws.Range("AI5").Value = 0.248
$ws.Range("AJ5").Value = 0.105
$ws.Range("AK5").Value = 0.324
$ws.Range("AU5").Value = 0.492
$ws.Range("AV5").Value = 0.08599999999999999
$ws.Range("AW5").Value = 0.294
$ws.Range("BA5").Value = 1.097
$ws.Range("BB5").Value = 0.01
$ws.Range("BC5").Value = 0.101
$ws.Range("BG5").Value = 0.356
$ws.Range("BH5").Value = 0.058
$ws.Range("BI5").Value = 0.241
$ws.Range("BM5").Value = 0.358
$ws.Range("BN5").Value = 0.014
$ws.Range("BO5").Value = 0.117
$ws.Range("BP5").Value = 0.366
$ws.Range("BQ5").Value = 0.409
$ws.Range("E5").Value = 0.602
$ws.Range("F5").Value = 0.065
$ws.Range("G5").Value = 0.254
$ws.Range("N5").Value = 0.585
$ws.Range("O5").Value = 0.04
$ws.Range("P5").Value = 0.2
$ws.Range("Q5").Value = 0.181
$ws.Range("S5").Value = 0.222
$ws.Range("W5").Value = 0.265
$ws.Range("X5").Value = 0.096
$ws.Range("Y5").Value = 0.309

# Row 6: recalculated stats
$ws.Range("AI6").Value = 0.254
$ws.Range("AU6").Value = 0.356
$ws.Range("BA6").Value = 1.428
$ws.Range("BG6").Value = 0.474
$ws.Range("BM6").Value = 0.476
$ws.Range("BP6").Value = 0.476
$ws.Range("BQ6").Value = 0.541
$ws.Range("E6").Value = 0.551
$ws.Range("N6").Value = 0.589
$ws.Range("Q6").Value = 0.218
$ws.Range("W6").Value = 0.318

# Row 7: recalculated stats
$ws.Range("AI7").Value = 0.25
$ws.Range("AU7").Value = 0.427
$ws.Range("BA7").Value = 1.208
$ws.Range("BG7").Value = 0.395
$ws.Range("BM7").Value = 0.397
$ws.Range("BP7").Value = 0.403
$ws.Range("BQ7").Value = 0.453
$ws.Range("E7").Value = 0.581
$ws.Range("N7").Value = 0.587
$ws.Range("Q7").Value = 0.194
$ws.Range("W7").Value = 0.284

# Row 8: recalculated stats
$ws.Range("AI8").Value = 0.319
$ws.Range("AJ8").Value = 0.159
$ws.Range("AK8").Value = 0.399
$ws.Range("AU8").Value = 0.508
$ws.Range("AV8").Value = 0.065
$ws.Range("AW8").Value = 0.255
$ws.Range("BA8").Value = 1.781
$ws.Range("BB8").Value = 0.052
$ws.Range("BC8").Value = 0.228
$ws.Range("BG8").Value = 0.573
$ws.Range("BH8").Value = 0.12
$ws.Range("BI8").Value = 0.346
$ws.Range("BM8").Value = 0.608
$ws.Range("BN8").Value = 0.062
$ws.Range("BO8").Value = 0.25
$ws.Range("BP8").Value = 0.594
$ws.Range("BQ8").Value = 0.625
$ws.Range("E8").Value = 0.713
$ws.Range("F8").Value = 0.083
$ws.Range("G8").Value = 0.288
$ws.Range("N8").Value = 0.8100000000000001
$ws.Range("O8").Value = 0.019
$ws.Range("P8").Value = 0.139
$ws.Range("Q8").Value = 0.199
$ws.Range("R8").Value = 0.097
$ws.Range("S8").Value = 0.312
$ws.Range("W8").Value = 0.459
$ws.Range("X8").Value = 0.127
$ws.Range("Y8").Value = 0.357

# Row 9: recalculated stats
$ws.Range("AI9").Value = 0.375
$ws.Range("AJ9").Value = 0.234
$ws.Range("AK9").Value = 0.484
$ws.Range("BA9").Value = 2
$ws.Range("BM9").Value = 0.75
$ws.Range("BN9").Value = 0.188
$ws.Range("BO9").Value = 0.433
$ws.Range("BP9").Value = 0.667
$ws.Range("BQ9").Value = 0.6830000000000001
$ws.Range("E9").Value = 0.75
$ws.Range("F9").Value = 0.188
$ws.Range("G9").Value = 0.433
$ws.Range("N9").Value = 0.875
$ws.Range("O9").Value = 0.109
$ws.Range("P9").Value = 0.331

# Row 10: recalculated stats
$ws.Range("AI10").Value = 0.375
$ws.Range("AJ10").Value = 0.234
$ws.Range("AK10").Value = 0.484
$ws.Range("BA10").Value = 2.25
$ws.Range("BB10").Value = 0.188
$ws.Range("BC10").Value = 0.433
$ws.Range("BM10").Value = 0.75
$ws.Range("BN10").Value = 0.188
$ws.Range("BO10").Value = 0.433
$ws.Range("BP10").Value = 0.75
$ws.Range("BQ10").Value = 0.778
$ws.Range("E10").Value = 0.75
$ws.Range("F10").Value = 0.188
$ws.Range("G10").Value = 0.433
$ws.Range("N10").Value = 1
$ws.Range("O10").Value = 0
$ws.Range("P10").Value = 0
$ws.Range("W10").Value = 0.625
$ws.Range("X10").Value = 0.234
$ws.Range("Y10").Value = 0.484

# Row 11: recalculated stats
$ws.Range("AI11").Value = 0.375
$ws.Range("AJ11").Value = 0.234
$ws.Range("AK11").Value = 0.484
$ws.Range("AU11").Value = 0.75
$ws.Range("AV11").Value = 0.188
$ws.Range("AW11").Value = 0.433
$ws.Range("BA11").Value = 2.25
$ws.Range("BB11").Value = 0.188
$ws.Range("BC11").Value = 0.433
$ws.Range("BM11").Value = 0.75
$ws.Range("BN11").Value = 0.188
$ws.Range("BO11").Value = 0.433
$ws.Range("BP11").Value = 0.75
$ws.Range("BQ11").Value = 0.778
$ws.Range("E11").Value = 0.875
$ws.Range("F11").Value = 0.109
$ws.Range("G11").Value = 0.331
$ws.Range("N11").Value = 1
$ws.Range("O11").Value = 0
$ws.Range("P11").Value = 0
$ws.Range("W11").Value = 0.625
$ws.Range("X11").Value = 0.234
$ws.Range("Y11").Value = 0.484

# Row 12: recalculated stats
$ws.Range("AU12").Value = 3.429
$ws.Range("AV12").Value = 4.245
$ws.Range("AW12").Value = 2.06
$ws.Range("BA12").Value = 3.5
$ws.Range("BB12").Value = 0.139
$ws.Range("BC12").Value = 0.373
$ws.Range("BP12").Value = 1.167
$ws.Range("BQ12").Value = 1.163
$ws.Range("E12").Value = 1.571
$ws.Range("F12").Value = 1.959
$ws.Range("G12").Value = 1.4
$ws.Range("N12").Value = 1.125
$ws.Range("O12").Value = 0.109
$ws.Range("P12").Value = 0.331
$ws.Range("W12").Value = 1.2
$ws.Range("X12").Value = 0.16
$ws.Range("Y12").Value = 0.4

# Row 13: recalculated stats
$ws.Range("BP13").Value = 0.646
$ws.Range("BQ13").Value = 0.553
